# templateEstudiantesMaes.xlsx — "preparado por fin un funcional de creacion
# de excel y pdf": the generator sheet is reset back to a blank template —
# the label cells in column B (CODIGO, NOMBRE, CORREO, ...) are cleared so
# only the title row remains in the shared-string table, the view is
# scrolled down to where the user was working, and the print setup is
# switched to landscape with narrower/asymmetric margins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Blank out the template's row labels (B10:B23) -------------------------
# ClearContents drops the shared-string value from each cell but keeps the
# cell's style (s="7"), matching the template-reset in the diff.
$ws.Range("B10:B23").ClearContents() | Out-Null

# --- Sheet view: scrolled to row 6, active cell D15 -------------------------
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D15").Select() | Out-Null

# --- Page setup: landscape, narrower/asymmetric margins --------------------
$ps = $ws.PageSetup
$ps.LeftMargin = 17.00787401574803
$ps.RightMargin = 413.8582677165354
$ps.TopMargin = 53.85826771653544
$ps.BottomMargin = 53.85826771653544
$ps.HeaderMargin = 22.677165354330707
$ps.FooterMargin = 22.677165354330707
$ps.Orientation = 2
